# refactor: ReadXlsx function name
#
# The workbook is a price-scraper log ("Vasculhador_de_Precos"). This run
# replaces the stale duplicated rows (the previous scrape re-saved the same
# Amazon/Kabum/Magazine Luiza rows twice) with freshly scraped rows for new
# products, while keeping each store's block of rows together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string (e.g. "5", "4.7", "5.0") to be
# stored as text rather than being auto-converted to a number, without
# leaving the cell's number format / style changed afterwards.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 4 - Amazon: new CPU
$ws.Cells.Item(4, 1).Value = "Processador AMD Ryzen 5 5600G, 3.9GHz (4.4GHz Max Turbo), AM4, Vídeo Integrado"
$ws.Cells.Item(4, 2).Value = 1099
$ws.Cells.Item(4, 3).Value = "4,8"
$ws.Cells.Item(4, 4).Value = "Amazon"
$ws.Cells.Item(4, 5).Value = "23/10/2024"

# Row 5 - Amazon: new motherboard
$ws.Cells.Item(5, 1).Value = "Placa Mãe Gigabyte B760M AORUS ELITE (rev. 1.0), LGA 1700, DDR5"
$ws.Cells.Item(5, 2).Value = 1078
$ws.Cells.Item(5, 3).Value = "4,4"
$ws.Cells.Item(5, 4).Value = "Amazon"
$ws.Cells.Item(5, 5).Value = "23/10/2024"

# Row 6 - Kabum: Console Playstation 5 (moved down from old row 4)
$ws.Cells.Item(6, 1).Value = "Console Playstation 5 Sony Slim, SSD 1TB, Controle Sem Fio Dualsense, Edição Digital, Branco, Returnal E Ratchet E Clank"
$ws.Cells.Item(6, 2).Value = 3799.04
Set-TextValue $ws.Cells.Item(6, 3) "5"
$ws.Cells.Item(6, 4).Value = "Kabum"
$ws.Cells.Item(6, 5).Value = "23/10/2024"

# Row 7 - Kabum: Placa de Vídeo RTX 3060 (review re-grouped to "5")
$ws.Cells.Item(7, 1).Value = "Placa de Vídeo RTX 3060 1-Click OC Galax NVIDIA GeForce, 12GB GDDR6, LHR, DLSS, Ray Tracing - 36NOL7MD1VOC"
$ws.Cells.Item(7, 2).Value = 1639.99
Set-TextValue $ws.Cells.Item(7, 3) "5"
$ws.Cells.Item(7, 4).Value = "Kabum"
$ws.Cells.Item(7, 5).Value = "23/10/2024"

# Row 8 - Kabum: new CPU
$ws.Cells.Item(8, 1).Value = "Processador AMD Ryzen 5 5600, 3.5GHz (4.4GHz Max Turbo), Cache 35MB, AM4, Sem Vídeo - 100-100000927BOX"
$ws.Cells.Item(8, 2).Value = 799.99
Set-TextValue $ws.Cells.Item(8, 3) "5"
$ws.Cells.Item(8, 4).Value = "Kabum"
$ws.Cells.Item(8, 5).Value = "23/10/2024"

# Row 9 - Kabum: new motherboard
$ws.Cells.Item(9, 1).Value = "Placa-Mãe AsRock B550M Steel Legend, AMD AM4 B550, DDR4 4733 OC, USB 3.2"
$ws.Cells.Item(9, 2).Value = 1612.53
Set-TextValue $ws.Cells.Item(9, 3) "5"
$ws.Cells.Item(9, 4).Value = "Kabum"
$ws.Cells.Item(9, 5).Value = "23/10/2024"

# Row 10 - Magazine Luiza: n/a placeholder row (moved down from old row 6)
$ws.Cells.Item(10, 1).Value = "n/a"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = "n/a"
$ws.Cells.Item(10, 4).Value = "Magazine Luiza"
$ws.Cells.Item(10, 5).Value = "23/10/2024"

# Row 11 - Magazine Luiza: Playstation 5 bundle (moved up from old row 13)
$ws.Cells.Item(11, 1).Value = "Playstation 5 Slim Digital Bundle Returnal + Ratchet & Clank"
$ws.Cells.Item(11, 2).Value = 3922.05
Set-TextValue $ws.Cells.Item(11, 3) "4.7"
$ws.Cells.Item(11, 4).Value = "Magazine Luiza"
$ws.Cells.Item(11, 5).Value = "23/10/2024"

# Row 12 - Magazine Luiza: new CPU
$ws.Cells.Item(12, 1).Value = "Processador AMD Ryzen 5 7600, 5.1GHz Max Turbo, Cache 38MB, AM5, 6 Núcleos, Vídeo Integrado - 100-100001015BOX"
$ws.Cells.Item(12, 2).Value = 1448.99
Set-TextValue $ws.Cells.Item(12, 3) "5.0"
$ws.Cells.Item(12, 4).Value = "Magazine Luiza"
$ws.Cells.Item(12, 5).Value = "23/10/2024"

# Row 13 - Magazine Luiza: n/a placeholder row, keeping the 4.7 review value
$ws.Cells.Item(13, 1).Value = "n/a"
$ws.Cells.Item(13, 2).Value = 0
Set-TextValue $ws.Cells.Item(13, 3) "4.7"
$ws.Cells.Item(13, 4).Value = "Magazine Luiza"
$ws.Cells.Item(13, 5).Value = "23/10/2024"
